$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 897; existing rows 897-936 shift down to 899-938.
$ws.Rows("897:898").Insert()

# New row 897 (Calidad: Primera)
$ws.Cells.Item(897, 1).Value = 3
$ws.Cells.Item(897, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(897, 3).Value = "Coquimbo"
$ws.Cells.Item(897, 4).Value = 44939
$ws.Cells.Item(897, 5).Value = 5
$ws.Cells.Item(897, 6).Value = 100112023
$ws.Cells.Item(897, 7).Value = "Brócoli"
$ws.Cells.Item(897, 8).Value = "Sin especificar"
$ws.Cells.Item(897, 9).Value = "Primera"
$ws.Cells.Item(897, 10).Value = 2900
$ws.Cells.Item(897, 11).Value = 800
$ws.Cells.Item(897, 12).Value = 850
$ws.Cells.Item(897, 13).Value = 819
$ws.Cells.Item(897, 14).Value = "`$/unidad"
$ws.Cells.Item(897, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(897, 16).Value = 819
$ws.Cells.Item(897, 17).Value = 1
$ws.Cells.Item(897, 18).Value = "Hortaliza"

# New row 898 (Calidad: Segunda)
$ws.Cells.Item(898, 1).Value = 3
$ws.Cells.Item(898, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(898, 3).Value = "Coquimbo"
$ws.Cells.Item(898, 4).Value = 44939
$ws.Cells.Item(898, 5).Value = 5
$ws.Cells.Item(898, 6).Value = 100112023
$ws.Cells.Item(898, 7).Value = "Brócoli"
$ws.Cells.Item(898, 8).Value = "Sin especificar"
$ws.Cells.Item(898, 9).Value = "Segunda"
$ws.Cells.Item(898, 10).Value = 1200
$ws.Cells.Item(898, 11).Value = 700
$ws.Cells.Item(898, 12).Value = 700
$ws.Cells.Item(898, 13).Value = 700
$ws.Cells.Item(898, 14).Value = "`$/unidad"
$ws.Cells.Item(898, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(898, 16).Value = 700
$ws.Cells.Item(898, 17).Value = 1
$ws.Cells.Item(898, 18).Value = "Hortaliza"
